# ---------------------------------------------------------------------------
# Split the single "Hoja1" table (works of art + their newspaper references)
# into two dedicated sheets: "Obra" (artwork) and "Referente" (press
# reference).  Hoja1 itself is left with its original data, it just stops
# being the active tab.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# xlPasteValues / xlPasteFormats constants
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# --- create the new sheets, right after Hoja1, in the desired order -------
$wsObra = $wb.Worksheets.Add($null, $ws1)
$wsObra.Name = "Obra"

$wsReferente = $wb.Worksheets.Add($null, $wsObra)
$wsReferente.Name = "Referente"

# ===========================================================================
# "Obra" sheet: ID, Titulo, Fecha, Dimensiones, Tecnica, Periodico, Archivo
# ===========================================================================

$wsObra.Columns.Item(2).ColumnWidth = (43.3828125 - 0.8333333333333334)
$wsObra.Columns.Item(3).ColumnWidth = (17.3046875 - 0.8333333333333334)
$wsObra.Columns.Item(4).ColumnWidth = (14.53515625 - 0.8333333333333334)
$wsObra.Columns.Item(5).ColumnWidth = (20 - 0.8333333333333334)
$wsObra.Columns.Item(6).ColumnWidth = (14.3828125 - 0.8333333333333334)
$wsObra.Columns.Item(7).ColumnWidth = (25.23046875 - 0.8333333333333334)

# header row
$wsObra.Range("A1").Value = "ID"
$wsObra.Range("B1").Value = "Título"
$wsObra.Range("C1").Value = "Fecha"
$wsObra.Range("D1").Value = "Dimensiones"
$wsObra.Range("E1").Value = "Técnica"
$wsObra.Range("F1").Value = "Periódico"
$wsObra.Range("G1").Value = "Archivo"

# row 2  <- Hoja1 row 2 (ID 1)
$wsObra.Range("A2").Value = 1
$wsObra.Range("B2").Value = "Los Suicidas del Sisga No 1"
$wsObra.Range("C2").Value = 1965
$wsObra.Range("D2").Value = "120x100cm"
$wsObra.Range("E2").Value = "Óleo sobre lienzo"
$wsObra.Range("G2").Value = "suicidas-sisga-1.jpg"
$ws1.Range("C2:G2").Copy()
$wsObra.Range("C2:G2").PasteSpecial($xlPasteFormats)
$ws1.Range("B2").Copy()
$wsObra.Range("B2").PasteSpecial($xlPasteFormats)
$wsObra.Range("B2").VerticalAlignment = -4107

# row 3  <- Hoja1 row 4 (ID 3)
$wsObra.Range("A3").Value = 3
$wsObra.Range("B3").Value = "El Paraíso"
$wsObra.Range("C3").Value = 1997
$wsObra.Range("D3").Value = "160x45 cm"
$wsObra.Range("E3").Value = "Óleo sobre lienzo"
$wsObra.Range("G3").Value = "el-paraiso.jpg"
$ws1.Range("C4:G4").Copy()
$wsObra.Range("C3:G3").PasteSpecial($xlPasteFormats)

# row 4  <- Hoja1 row 7 (ID 6)
$wsObra.Range("A4").Value = 6
$wsObra.Range("B4").Value = "Zócalo de la tragedia"
$wsObra.Range("C4").Value = 1983
$wsObra.Range("D4").Value = "100x70"
$wsObra.Range("E4").Value = "Tipografía sobre papel"
$wsObra.Range("G4").Value = "zocalo-tragedia.jpg"
$ws1.Range("C7:G7").Copy()
$wsObra.Range("C4:G4").PasteSpecial($xlPasteFormats)

$wsObra.Range("B3").Select()

# ===========================================================================
# "Referente" sheet: ID, Titulo, Fecha, Periodico, Archivo
# ===========================================================================

$wsReferente.Columns.Item(2).ColumnWidth = (43.3828125 - 0.8333333333333334)
$wsReferente.Columns.Item(3).ColumnWidth = (17.3046875 - 0.8333333333333334)
$wsReferente.Columns.Item(4).ColumnWidth = (14.3828125 - 0.8333333333333334)
$wsReferente.Columns.Item(5).ColumnWidth = (25.23046875 - 0.8333333333333334)

# header row
$wsReferente.Range("A1").Value = "ID"
$wsReferente.Range("B1").Value = "Título"
$wsReferente.Range("C1").Value = "Fecha"
$wsReferente.Range("D1").Value = "Periódico"
$wsReferente.Range("E1").Value = "Archivo"

# row 2 <- Hoja1 row 3 (ID 2)
$wsReferente.Range("A2").Value = 2
$wsReferente.Range("B2").Value = "Doble suicidio en ""El Sisga"""
$wsReferente.Range("C2").Value = "Junio 29 1965"
$wsReferente.Range("D2").Value = "El Tiempo"
$wsReferente.Range("E2").Value = "doble-suicidio-el-tiempo.jpg"
$ws1.Range("C3").Copy()
$wsReferente.Range("C2").PasteSpecial($xlPasteFormats)
$ws1.Range("F3:G3").Copy()
$wsReferente.Range("D2:E2").PasteSpecial($xlPasteFormats)

# row 3 <- Hoja1 row 5 (ID 4)
$wsReferente.Range("A3").Value = 4
$wsReferente.Range("B3").Value = "Una indígena y su hijo murieron en persecución"
$wsReferente.Range("C3").Value = "Mayo 24 del 96"
$wsReferente.Range("D3").Value = "El Tiempo"
$wsReferente.Range("E3").Value = "indigena-hijo-el-tiempo.jpg"
$ws1.Range("C5").Copy()
$wsReferente.Range("C3").PasteSpecial($xlPasteFormats)
$ws1.Range("F5:G5").Copy()
$wsReferente.Range("D3:E3").PasteSpecial($xlPasteFormats)

# row 4 <- Hoja1 row 6 (ID 5)
$wsReferente.Range("A4").Value = 5
$wsReferente.Range("B4").Value = "Láminas de paisajes latinoamericanos"
$wsReferente.Range("E4").Value = "laminas-paisajes.jpg"
$ws1.Range("C6").Copy()
$wsReferente.Range("C4").PasteSpecial($xlPasteFormats)
$ws1.Range("F6:G6").Copy()
$wsReferente.Range("D4:E4").PasteSpecial($xlPasteFormats)

# row 5 <- Hoja1 row 8 (ID 7)
$wsReferente.Range("A5").Value = 7
$wsReferente.Range("B5").Value = "Exmilitar Mata a la Esposa de su Amigo y se Suicida"
$wsReferente.Range("E5").Value = "exmilitar-mata-esposa.jpg"

$wsReferente.Range("C11").Select()

# Referente becomes the active tab, Hoja1 automatically loses tabSelected.
$wsReferente.Activate()
